$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a1"
$ws.Range("C2").Value = "Itga2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 6.072131
$ws.Range("H2").Value = 18.216393
$ws.Range("I2").Value = 0.003943999267036455
$ws.Range("J2").Value = 0.003943999267036454
$ws.Range("K2").Value = 2.0
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.242496666666667
$ws.Range("N2").Value = 3.72749
$ws.Range("O2").Value = 0.3272238221337332
$ws.Range("P2").Value = 0.3272238221337332
$ws.Range("Q2").Value = 7.544602527063333
$ws.Range("R2").Value = 67.90142274357
$ws.Range("S2").Value = 0.001290570514652311
$ws.Range("T2").Value = 0.001290570514652311

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a1"
$ws.Range("C3").Value = "Itga2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 6.072131
$ws.Range("H3").Value = 18.216393
$ws.Range("I3").Value = 0.003943999267036455
$ws.Range("J3").Value = 0.003943999267036454
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 1.744414
$ws.Range("N3").Value = 5.233242000000001
$ws.Range("O3").Value = 0.459408730644692
$ws.Range("P3").Value = 0.459408730644692
$ws.Range("Q3").Value = 10.592310326234
$ws.Range("R3").Value = 95.33079293610601
$ws.Range("S3").Value = 0.001811907696932813
$ws.Range("T3").Value = 0.001811907696932813

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a1"
$ws.Range("C4").Value = "Itga2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 6.072131
$ws.Range("H4").Value = 18.216393
$ws.Range("I4").Value = 0.003943999267036455
$ws.Range("J4").Value = 0.003943999267036454
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 0.8101743333333333
$ws.Range("N4").Value = 2.430523
$ws.Range("O4").Value = 0.2133674472215748
$ws.Range("P4").Value = 0.2133674472215748
$ws.Range("Q4").Value = 4.919484684837666
$ws.Range("R4").Value = 44.275362163539
$ws.Range("S4").Value = 0.0008415210554513304
$ws.Range("T4").Value = 0.0008415210554513302

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col1a1"
$ws.Range("C5").Value = "Itga2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 1480.851806666667
$ws.Range("H5").Value = 4442.55542
$ws.Range("I5").Value = 0.9618498744646554
$ws.Range("J5").Value = 0.9618498744646552
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.242496666666667
$ws.Range("N5").Value = 3.72749
$ws.Range("O5").Value = 0.3272238221337332
$ws.Range("P5").Value = 0.3272238221337332
$ws.Range("Q5").Value = 1839.953433610644
$ws.Range("R5").Value = 16559.5809024958
$ws.Range("S5").Value = 0.3147401922411761
$ws.Range("T5").Value = 0.314740192241176

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a1"
$ws.Range("C6").Value = "Itga2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 1480.851806666667
$ws.Range("H6").Value = 4442.55542
$ws.Range("I6").Value = 0.9618498744646554
$ws.Range("J6").Value = 0.9618498744646552
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 1.744414
$ws.Range("N6").Value = 5.233242000000001
$ws.Range("O6").Value = 0.459408730644692
$ws.Range("P6").Value = 0.459408730644692
$ws.Range("Q6").Value = 2583.218623474627
$ws.Range("R6").Value = 23248.96761127164
$ws.Range("S6").Value = 0.4418822298985636
$ws.Range("T6").Value = 0.4418822298985636

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a1"
$ws.Range("C7").Value = "Itga2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 1480.851806666667
$ws.Range("H7").Value = 4442.55542
$ws.Range("I7").Value = 0.9618498744646554
$ws.Range("J7").Value = 0.9618498744646552
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 0.8101743333333333
$ws.Range("N7").Value = 2.430523
$ws.Range("O7").Value = 0.2133674472215748
$ws.Range("P7").Value = 0.2133674472215748
$ws.Range("Q7").Value = 1199.748125231629
$ws.Range("R7").Value = 10797.73312708466
$ws.Range("S7").Value = 0.2052274523249157
$ws.Range("T7").Value = 0.2052274523249156

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col1a1"
$ws.Range("C8").Value = "Itga2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 52.663316
$ws.Range("H8").Value = 157.989948
$ws.Range("I8").Value = 0.03420612626830831
$ws.Range("J8").Value = 0.0342061262683083
$ws.Range("K8").Value = 2.0
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.242496666666667
$ws.Range("N8").Value = 3.72749
$ws.Range("O8").Value = 0.3272238221337332
$ws.Range("P8").Value = 0.3272238221337332
$ws.Range("Q8").Value = 65.43399458561333
$ws.Range("R8").Value = 588.90595127052
$ws.Range("S8").Value = 0.01119305937790494
$ws.Range("T8").Value = 0.01119305937790494

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col1a1"
$ws.Range("C9").Value = "Itga2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 52.663316
$ws.Range("H9").Value = 157.989948
$ws.Range("I9").Value = 0.03420612626830831
$ws.Range("J9").Value = 0.0342061262683083
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 1.744414
$ws.Range("N9").Value = 5.233242000000001
$ws.Range("O9").Value = 0.459408730644692
$ws.Range("P9").Value = 0.459408730644692
$ws.Range("Q9").Value = 91.866625716824
$ws.Range("R9").Value = 826.7996314514161
$ws.Range("S9").Value = 0.01571459304919557
$ws.Range("T9").Value = 0.01571459304919557

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col1a1"
$ws.Range("C10").Value = "Itga2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 52.663316
$ws.Range("H10").Value = 157.989948
$ws.Range("I10").Value = 0.03420612626830831
$ws.Range("J10").Value = 0.0342061262683083
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 0.8101743333333333
$ws.Range("N10").Value = 2.430523
$ws.Range("O10").Value = 0.2133674472215748
$ws.Range("P10").Value = 0.2133674472215748
$ws.Range("Q10").Value = 42.66646693142267
$ws.Range("R10").Value = 383.998202382804
$ws.Range("S10").Value = 0.007298473841207797
$ws.Range("T10").Value = 0.007298473841207794

Write-Output "Applied Col1a1-Itga2 updates (Dr Hou advice)"
